$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New master data rows (regcntr_id, machine_id, lang_code, is_active, cr_by, cr_dtimes)
# regcntr_id restarts at 10002 while machine_id continues sequentially from 10021.
$newRows = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $pair = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Keep the sheet tab selected and move the active selection / view down to
# the first empty row below the data, as the author's last on-screen state.
$ws.Range("A31:XFD1048576").Select() | Out-Null

# Switch the page to portrait orientation (as set via Page Setup).
$ws.PageSetup.Orientation = 1
